# Update cryptos list: set Price (D) and Volume(1h) (E) columns with new values.
# All target cells hold inline strings, so we force text format to preserve exact
# formatting (leading/trailing spaces, percent signs, dotted "thousands" style prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '21.797.57'
    'E2' = '  -1.62%  '
    'D3' = '1.543.20'
    'E3' = '  -1.14%  '
    'D4' = '0.9993'
    'E4' = '  -0.14%  '
    'E5' = '  -0.07%  '
    'D6' = '289.96'
    'E6' = '  -0.04%  '
    'D7' = '0.3903'
    'E7' = '  +2.66%  '
    'D8' = '0.3190'
    'E8' = '  -2.84%  '
    'D9' = '43.27'
    'E9' = '  -1.05%  '
    'D10' = '0.07207'
    'E10' = '  -2.14%  '
    'D11' = '1.069'
    'E11' = '  -6.12%  '
    'D12' = '0.9993'
    'D13' = '5.638'
    'D14' = '18.66'
    'E14' = '  -6.50%  '
    'D15' = '6.610'
    'E15' = '  -3.79%  '
    'D16' = '1.539.56'
    'E16' = '  -1.48%  '
    'D17' = '0.00001110'
    'E17' = '  +1.43%  '
    'D18' = '0.06583'
    'E18' = '  -0.80%  '
    'D19' = '83.37'
    'E19' = '  -2.46%  '
    'D20' = '0.9998'
    'E20' = '  -0.08%  '
    'D21' = '6.161'
    'E21' = '  -4.52%  '
    'D22' = '15.41'
    'E22' = '  -4.42%  '
    'D23' = '10.89'
    'E23' = '  -7.27%  '
    'D24' = '2.367'
    'E24' = '  +4.47%  '
    'D25' = '21.796.09'
    'E25' = '  -1.63%  '
    'D26' = '2.399'
    'E26' = '  -5.43%  '
    'D27' = '145.34'
    'E27' = '  -3.78%  '
    'D28' = '18.43'
    'E28' = '  -3.37%  '
    'D29' = '4.852'
    'E29' = '  -0.38%  '
    'D30' = '1.716.36'
    'E30' = '  -1.29%  '
    'D31' = '117.85'
    'E31' = '  -2.91%  '
    'D32' = '0.9712'
    'E32' = '  -13.25%  '
    'D33' = '5.922'
    'E33' = '  -1.87%  '
    'D34' = '0.08203'
    'E34' = '  -0.08%  '
    'D35' = '9.005'
    'E35' = '  -3.64%  '
    'D36' = '0.06134'
    'E36' = '  -1.36%  '
    'D37' = '5.136'
    'E37' = '  -2.89%  '
    'E38' = '  -3.98%  '
    'D39' = '0.2045'
    'E39' = '  -4.38%  '
    'D40' = '1.185'
    'E40' = '  -3.78%  '
    'D41' = '1.429'
    'E41' = '  -24.13%  '
    'D42' = '0.9997'
    'E42' = '  -0.10%  '
    'D43' = '10.69'
    'E43' = '  -3.44%  '
    'D44' = '0.5779'
    'E44' = '  -3.40%  '
    'D45' = '13.10'
    'E45' = '  -4.35%  '
    'D46' = '3.739'
    'E46' = '  -0.50%  '
    'D47' = '0.5544'
    'E47' = '  -4.25%  '
    'D48' = '117.96'
    'E48' = '  -2.49%  '
    'D49' = '1.882'
    'E49' = '  -5.34%  '
    'E50' = '  -3.26%  '
    'D51' = '0.06737'
    'E51' = '  -3.65%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
